$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1) Update the CO number on the existing schedule line (row 2)
$ws.Range("A2").Value = "3013992465"

# 2) Add a new schedule line (row 3) mirroring row 2's data / formatting
$ws.Range("A3").Value = "3013997547"
$ws.Range("B3").Value = "TA5TBBPC"
$ws.Range("C3").Value = "102717"
$ws.Range("D3").Value = "1"
$ws.Range("E3").Value = "EA"
$ws.Range("F3").Value = "20-Firm"
$ws.Range("G3").Value = "102717"

# Match row 3's number formatting / alignment to row 2's (C:G), same as the
# existing schedule line above it
$ws.Range("C3:G3").NumberFormat = "@"
$ws.Range("C3").HorizontalAlignment = -4131
$ws.Range("C3").VerticalAlignment = -4160
$ws.Range("D3:E3").HorizontalAlignment = -4152
$ws.Range("D3:E3").VerticalAlignment = -4160
$ws.Range("F3:G3").HorizontalAlignment = -4131
$ws.Range("F3:G3").VerticalAlignment = -4160

# 3) Select just the updated cell, and autofit column A to its new content
$ws.Range("A2").Select()
$ws.Columns.Item(1).AutoFit()
